# Generate Report for Handoff
# Refresh the "Latest Handoff Datetime" value for the 5ead2ef3... file row
# (row 5, "Ready for handoff") on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-10 16:36:07"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-10 16:36:15"
